# Update PLC data 2025-10-13 14:01:20
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C2").Value = 275
$ws.Range("C3").Value = 172265
$ws.Range("C4").Value = 163056
$ws.Range("C8").Value = 65.88
